$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 15 <-> row 16
$ws.Range("B15").Value = 7004588
$ws.Range("B16").Value = 7004589
$ws.Range("F15").Value = "Umm Salal"
$ws.Range("F16").Value = "AlMuaidar"
$ws.Range("G15").Value = "Qatar SC Doha"
$ws.Range("G16").Value = "Al Markhiya"
$ws.Range("H15").Value = 2
$ws.Range("H16").Value = 5
$ws.Range("I15").Value = 2
$ws.Range("I16").Value = 2
$ws.Range("J15").Value = "D"
$ws.Range("J16").Value = "H"
$ws.Range("K15").Value = 2.9
$ws.Range("K16").Value = 2.4
$ws.Range("L15").Value = 3.3
$ws.Range("L16").Value = 3.6
$ws.Range("M15").Value = 2.25
$ws.Range("M16").Value = 2.375
$ws.Range("N15").Value = 3.75
$ws.Range("N16").Value = 2.5
$ws.Range("O15").Value = 3.4
$ws.Range("O16").Value = 3.5
$ws.Range("P15").Value = 1.909
$ws.Range("P16").Value = 2.3
$ws.Range("Q15").Value = 0.5
$ws.Range("Q16").Value = 0
$ws.Range("R15").Value = 1.85
$ws.Range("R16").Value = 1.975
$ws.Range("S15").Value = 1.95
$ws.Range("S16").Value = 1.825
$ws.Range("T15").Value = 2.5
$ws.Range("T16").Value = 3
$ws.Range("U15").Value = 1.85
$ws.Range("U16").Value = 2
$ws.Range("V15").Value = 1.95
$ws.Range("V16").Value = 1.8
$ws.Range("W15").Value = -1
$ws.Range("W16").Value = 1.5
$ws.Range("X15").Value = 2.4
$ws.Range("X16").Value = -1
$ws.Range("Y15").Value = -1
$ws.Range("Y16").Value = -1
$ws.Range("Z15").Value = 0.8500000000000001
$ws.Range("Z16").Value = 0.9750000000000001
$ws.Range("AA15").Value = -1
$ws.Range("AA16").Value = -1
$ws.Range("AB15").Value = 0.8500000000000001
$ws.Range("AB16").Value = 1
$ws.Range("AC15").Value = -1
$ws.Range("AC16").Value = -1

# Swap row 18 <-> row 19
$ws.Range("B18").Value = 7003585
$ws.Range("B19").Value = 7004591
$ws.Range("F18").Value = "Al Sadd"
$ws.Range("F19").Value = "AlShamal SC"
$ws.Range("G18").Value = "AlWakrah SC"
$ws.Range("G19").Value = "AlRayyan SC"
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 3
$ws.Range("I18").Value = 0
$ws.Range("I19").Value = 4
$ws.Range("J18").Value = "D"
$ws.Range("J19").Value = "A"
$ws.Range("K18").Value = 1.615
$ws.Range("K19").Value = 4.5
$ws.Range("L18").Value = 4
$ws.Range("L19").Value = 4.2
$ws.Range("M18").Value = 4.333
$ws.Range("M19").Value = 1.55
$ws.Range("N18").Value = 1.533
$ws.Range("N19").Value = 3.3
$ws.Range("O18").Value = 4.2
$ws.Range("O19").Value = 3.8
$ws.Range("P18").Value = 5
$ws.Range("P19").Value = 1.85
$ws.Range("Q18").Value = -1
$ws.Range("Q19").Value = 0.5
$ws.Range("R18").Value = 1.8
$ws.Range("R19").Value = 1.85
$ws.Range("S18").Value = 2
$ws.Range("S19").Value = 1.95
$ws.Range("T18").Value = 3.5
$ws.Range("T19").Value = 2.75
$ws.Range("U18").Value = 1.925
$ws.Range("U19").Value = 1.85
$ws.Range("V18").Value = 1.875
$ws.Range("V19").Value = 1.95
$ws.Range("W18").Value = -1
$ws.Range("W19").Value = -1
$ws.Range("X18").Value = 3.2
$ws.Range("X19").Value = -1
$ws.Range("Y18").Value = -1
$ws.Range("Y19").Value = 0.8500000000000001
$ws.Range("Z18").Value = -1
$ws.Range("Z19").Value = -1
$ws.Range("AA18").Value = 1
$ws.Range("AA19").Value = 0.95
$ws.Range("AB18").Value = -1
$ws.Range("AB19").Value = 0.8500000000000001
$ws.Range("AC18").Value = 0.875
$ws.Range("AC19").Value = -1

# Swap row 27 <-> row 28
$ws.Range("B27").Value = 7004597
$ws.Range("B28").Value = 7004596
$ws.Range("F27").Value = "AlArabi Doha"
$ws.Range("F28").Value = "Al Markhiya"
$ws.Range("G27").Value = "Qatar SC Doha"
$ws.Range("G28").Value = "Al Duhail"
$ws.Range("H27").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("I27").Value = 0
$ws.Range("I28").Value = 2
$ws.Range("J27").Value = "H"
$ws.Range("J28").Value = "A"
$ws.Range("K27").Value = 1.833
$ws.Range("K28").Value = 6.5
$ws.Range("L27").Value = 4
$ws.Range("L28").Value = 6
$ws.Range("M27").Value = 3.3
$ws.Range("M28").Value = 1.3
$ws.Range("N27").Value = 1.65
$ws.Range("N28").Value = 6.5
$ws.Range("O27").Value = 4.2
$ws.Range("O28").Value = 6
$ws.Range("P27").Value = 4
$ws.Range("P28").Value = 1.3
$ws.Range("Q27").Value = -0.75
$ws.Range("Q28").Value = 1.75
$ws.Range("R27").Value = 1.825
$ws.Range("R28").Value = 1.775
$ws.Range("S27").Value = 1.975
$ws.Range("S28").Value = 2.025
$ws.Range("T27").Value = 3
$ws.Range("T28").Value = 3.5
$ws.Range("U27").Value = 2
$ws.Range("U28").Value = 1.975
$ws.Range("V27").Value = 1.8
$ws.Range("V28").Value = 1.825
$ws.Range("W27").Value = 0.6499999999999999
$ws.Range("W28").Value = -1
$ws.Range("X27").Value = -1
$ws.Range("X28").Value = -1
$ws.Range("Y27").Value = -1
$ws.Range("Y28").Value = 0.3
$ws.Range("Z27").Value = 0.4125
$ws.Range("Z28").Value = 0.7749999999999999
$ws.Range("AA27").Value = -0.5
$ws.Range("AA28").Value = -1
$ws.Range("AB27").Value = -1
$ws.Range("AB28").Value = -1
$ws.Range("AC27").Value = 0.8
$ws.Range("AC28").Value = 0.825

# Swap row 40 <-> row 41
$ws.Range("B40").Value = 7004607
$ws.Range("B41").Value = 7004604
$ws.Range("F40").Value = "AlWakrah SC"
$ws.Range("F41").Value = "Qatar SC Doha"
$ws.Range("G40").Value = "Umm Salal"
$ws.Range("G41").Value = "Al Markhiya"
$ws.Range("H40").Value = 2
$ws.Range("H41").Value = 4
$ws.Range("I40").Value = 1
$ws.Range("I41").Value = 0
$ws.Range("J40").Value = "H"
$ws.Range("J41").Value = "H"
$ws.Range("K40").Value = 1.65
$ws.Range("K41").Value = 1.727
$ws.Range("L40").Value = 3.75
$ws.Range("L41").Value = 3.75
$ws.Range("M40").Value = 4.5
$ws.Range("M41").Value = 4
$ws.Range("N40").Value = 1.75
$ws.Range("N41").Value = 1.75
$ws.Range("O40").Value = 3.6
$ws.Range("O41").Value = 3.75
$ws.Range("P40").Value = 4
$ws.Range("P41").Value = 3.8
$ws.Range("Q40").Value = -0.75
$ws.Range("Q41").Value = -0.5
$ws.Range("R40").Value = 1.95
$ws.Range("R41").Value = 1.75
$ws.Range("S40").Value = 1.85
$ws.Range("S41").Value = 1.95
$ws.Range("T40").Value = 3
$ws.Range("T41").Value = 3
$ws.Range("U40").Value = 1.9
$ws.Range("U41").Value = 1.95
$ws.Range("V40").Value = 1.9
$ws.Range("V41").Value = 1.85
$ws.Range("W40").Value = 0.75
$ws.Range("W41").Value = 0.75
$ws.Range("X40").Value = -1
$ws.Range("X41").Value = -1
$ws.Range("Y40").Value = -1
$ws.Range("Y41").Value = -1
$ws.Range("Z40").Value = 0.475
$ws.Range("Z41").Value = 0.75
$ws.Range("AA40").Value = -0.5
$ws.Range("AA41").Value = -1
$ws.Range("AB40").Value = 0
$ws.Range("AB41").Value = 0.95
$ws.Range("AC40").Value = -0
$ws.Range("AC41").Value = -1

# Swap row 42 <-> row 43
$ws.Range("B42").Value = 7003590
$ws.Range("B43").Value = 7004611
$ws.Range("F42").Value = "Al Sadd"
$ws.Range("F43").Value = "AlMuaidar"
$ws.Range("G42").Value = "Al Markhiya"
$ws.Range("G43").Value = "Umm Salal"
$ws.Range("H42").Value = 5
$ws.Range("H43").Value = 1
$ws.Range("I42").Value = 0
$ws.Range("I43").Value = 3
$ws.Range("J42").Value = "H"
$ws.Range("J43").Value = "A"
$ws.Range("K42").Value = 1.125
$ws.Range("K43").Value = 4
$ws.Range("L42").Value = 8
$ws.Range("L43").Value = 3.25
$ws.Range("M42").Value = 10
$ws.Range("M43").Value = 1.8
$ws.Range("N42").Value = 1.2
$ws.Range("N43").Value = 3.75
$ws.Range("O42").Value = 5.75
$ws.Range("O43").Value = 3.4
$ws.Range("P42").Value = 10
$ws.Range("P43").Value = 1.8
$ws.Range("Q42").Value = -2
$ws.Range("Q43").Value = 0.5
$ws.Range("R42").Value = 1.975
$ws.Range("R43").Value = 1.925
$ws.Range("S42").Value = 1.825
$ws.Range("S43").Value = 1.875
$ws.Range("T42").Value = 3.5
$ws.Range("T43").Value = 3
$ws.Range("U42").Value = 1.95
$ws.Range("U43").Value = 1.85
$ws.Range("V42").Value = 1.85
$ws.Range("V43").Value = 1.95
$ws.Range("W42").Value = 0.2
$ws.Range("W43").Value = -1
$ws.Range("X42").Value = -1
$ws.Range("X43").Value = -1
$ws.Range("Y42").Value = -1
$ws.Range("Y43").Value = 0.8
$ws.Range("Z42").Value = 0.9750000000000001
$ws.Range("Z43").Value = -1
$ws.Range("AA42").Value = -1
$ws.Range("AA43").Value = 0.875
$ws.Range("AB42").Value = 0.95
$ws.Range("AB43").Value = 0.8500000000000001
$ws.Range("AC42").Value = -1
$ws.Range("AC43").Value = -1

# Swap row 56 <-> row 57
$ws.Range("B56").Value = 7004619
$ws.Range("B57").Value = 7004618
$ws.Range("F56").Value = "Qatar SC Doha"
$ws.Range("F57").Value = "Al Gharafa"
$ws.Range("G56").Value = "AlMuaidar"
$ws.Range("G57").Value = "AlWakrah SC"
$ws.Range("H56").Value = 3
$ws.Range("H57").Value = 1
$ws.Range("I56").Value = 2
$ws.Range("I57").Value = 1
$ws.Range("J56").Value = "H"
$ws.Range("J57").Value = "D"
$ws.Range("K56").Value = 1.4
$ws.Range("K57").Value = 2.5
$ws.Range("L56").Value = 4.5
$ws.Range("L57").Value = 3.75
$ws.Range("M56").Value = 5.5
$ws.Range("M57").Value = 2.4
$ws.Range("N56").Value = 1.833
$ws.Range("N57").Value = 1.95
$ws.Range("O56").Value = 3.8
$ws.Range("O57").Value = 3.8
$ws.Range("P56").Value = 3.25
$ws.Range("P57").Value = 3.25
$ws.Range("Q56").Value = -0.5
$ws.Range("Q57").Value = -0.5
$ws.Range("R56").Value = 1.9
$ws.Range("R57").Value = 1.975
$ws.Range("S56").Value = 1.9
$ws.Range("S57").Value = 1.825
$ws.Range("T56").Value = 3
$ws.Range("T57").Value = 3.5
$ws.Range("U56").Value = 2
$ws.Range("U57").Value = 1.975
$ws.Range("V56").Value = 1.8
$ws.Range("V57").Value = 1.825
$ws.Range("W56").Value = 0.833
$ws.Range("W57").Value = -1
$ws.Range("X56").Value = -1
$ws.Range("X57").Value = 2.8
$ws.Range("Y56").Value = -1
$ws.Range("Y57").Value = -1
$ws.Range("Z56").Value = 0.8999999999999999
$ws.Range("Z57").Value = -1
$ws.Range("AA56").Value = -1
$ws.Range("AA57").Value = 0.825
$ws.Range("AB56").Value = 1
$ws.Range("AB57").Value = -1
$ws.Range("AC56").Value = -1
$ws.Range("AC57").Value = 0.825

# Swap row 62 <-> row 63
$ws.Range("B62").Value = 7004622
$ws.Range("B63").Value = 7004621
$ws.Range("F62").Value = "Al Duhail"
$ws.Range("F63").Value = "Al Gharafa"
$ws.Range("G62").Value = "Qatar SC Doha"
$ws.Range("G63").Value = "AlShamal SC"
$ws.Range("H62").Value = 1
$ws.Range("H63").Value = 1
$ws.Range("I62").Value = 1
$ws.Range("I63").Value = 1
$ws.Range("J62").Value = "D"
$ws.Range("J63").Value = "D"
$ws.Range("K62").Value = 1.5
$ws.Range("K63").Value = 1.5
$ws.Range("L62").Value = 4.5
$ws.Range("L63").Value = 4.5
$ws.Range("M62").Value = 5
$ws.Range("M63").Value = 5
$ws.Range("N62").Value = 1.8
$ws.Range("N63").Value = 1.363
$ws.Range("O62").Value = 4
$ws.Range("O63").Value = 5
$ws.Range("P62").Value = 3.6
$ws.Range("P63").Value = 6
$ws.Range("Q62").Value = -0.5
$ws.Range("Q63").Value = -1.5
$ws.Range("R62").Value = 1.8
$ws.Range("R63").Value = 1.975
$ws.Range("S62").Value = 2
$ws.Range("S63").Value = 1.825
$ws.Range("T62").Value = 3
$ws.Range("T63").Value = 3.75
$ws.Range("U62").Value = 1.8
$ws.Range("U63").Value = 1.975
$ws.Range("V62").Value = 2
$ws.Range("V63").Value = 1.825
$ws.Range("W62").Value = -1
$ws.Range("W63").Value = -1
$ws.Range("X62").Value = 3
$ws.Range("X63").Value = 4
$ws.Range("Y62").Value = -1
$ws.Range("Y63").Value = -1
$ws.Range("Z62").Value = -1
$ws.Range("Z63").Value = -1
$ws.Range("AA62").Value = 1
$ws.Range("AA63").Value = 0.825
$ws.Range("AB62").Value = -1
$ws.Range("AB63").Value = -1
$ws.Range("AC62").Value = 1
$ws.Range("AC63").Value = 0.825

# Swap row 70 <-> row 71
$ws.Range("B70").Value = 7004627
$ws.Range("B71").Value = 7609335
$ws.Range("F70").Value = "AlRayyan SC"
$ws.Range("F71").Value = "AlShamal SC"
$ws.Range("G70").Value = "Al Markhiya"
$ws.Range("G71").Value = "AlArabi Doha"
$ws.Range("H70").Value = 6
$ws.Range("H71").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J70").Value = "H"
$ws.Range("J71").Value = "D"
$ws.Range("K70").Value = 1.444
$ws.Range("K71").Value = 4.5
$ws.Range("L70").Value = 4.75
$ws.Range("L71").Value = 4.2
$ws.Range("M70").Value = 5.75
$ws.Range("M71").Value = 1.533
$ws.Range("N70").Value = 1.363
$ws.Range("N71").Value = 5.25
$ws.Range("O70").Value = 5
$ws.Range("O71").Value = 4.2
$ws.Range("P70").Value = 6.5
$ws.Range("P71").Value = 1.45
$ws.Range("Q70").Value = -1.5
$ws.Range("Q71").Value = 1
$ws.Range("R70").Value = 2
$ws.Range("R71").Value = 2
$ws.Range("S70").Value = 1.8
$ws.Range("S71").Value = 1.8
$ws.Range("T70").Value = 3.25
$ws.Range("T71").Value = 3
$ws.Range("U70").Value = 1.95
$ws.Range("U71").Value = 1.9
$ws.Range("V70").Value = 1.85
$ws.Range("V71").Value = 1.9
$ws.Range("W70").Value = 0.363
$ws.Range("W71").Value = -1
$ws.Range("X70").Value = -1
$ws.Range("X71").Value = 3.2
$ws.Range("Y70").Value = -1
$ws.Range("Y71").Value = -1
$ws.Range("Z70").Value = 1
$ws.Range("Z71").Value = 1
$ws.Range("AA70").Value = -1
$ws.Range("AA71").Value = -1
$ws.Range("AB70").Value = 0.95
$ws.Range("AB71").Value = -1
$ws.Range("AC70").Value = -1
$ws.Range("AC71").Value = 0.8999999999999999

# Swap row 81 <-> row 82
$ws.Range("B81").Value = 7840802
$ws.Range("B82").Value = 7840805
$ws.Range("F81").Value = "AlAhli Doha"
$ws.Range("F82").Value = "AlWakrah SC"
$ws.Range("G81").Value = "Umm Salal"
$ws.Range("G82").Value = "Al Markhiya"
$ws.Range("H81").Value = 1
$ws.Range("H82").Value = 1
$ws.Range("I81").Value = 2
$ws.Range("I82").Value = 2
$ws.Range("J81").Value = "A"
$ws.Range("J82").Value = "A"
$ws.Range("K81").Value = 2.4
$ws.Range("K82").Value = 1.062
$ws.Range("L81").Value = 4
$ws.Range("L82").Value = 11
$ws.Range("M81").Value = 2.25
$ws.Range("M82").Value = 17
$ws.Range("N81").Value = 2.3
$ws.Range("N82").Value = 1.363
$ws.Range("O81").Value = 4
$ws.Range("O82").Value = 4.75
$ws.Range("P81").Value = 2.375
$ws.Range("P82").Value = 7
$ws.Range("Q81").Value = 0
$ws.Range("Q82").Value = -1.25
$ws.Range("R81").Value = 1.875
$ws.Range("R82").Value = 1.75
$ws.Range("S81").Value = 1.925
$ws.Range("S82").Value = 1.95
$ws.Range("T81").Value = 3.25
$ws.Range("T82").Value = 3.25
$ws.Range("U81").Value = 2
$ws.Range("U82").Value = 1.975
$ws.Range("V81").Value = 1.8
$ws.Range("V82").Value = 1.825
$ws.Range("W81").Value = -1
$ws.Range("W82").Value = -1
$ws.Range("X81").Value = -1
$ws.Range("X82").Value = -1
$ws.Range("Y81").Value = 1.375
$ws.Range("Y82").Value = 6
$ws.Range("Z81").Value = -1
$ws.Range("Z82").Value = -1
$ws.Range("AA81").Value = 0.925
$ws.Range("AA82").Value = 0.95
$ws.Range("AB81").Value = -0.5
$ws.Range("AB82").Value = -0.5
$ws.Range("AC81").Value = 0.4
$ws.Range("AC82").Value = 0.4125

# Swap row 90 <-> row 91
$ws.Range("B90").Value = 7840810
$ws.Range("B91").Value = 7840809
$ws.Range("F90").Value = "AlRayyan SC"
$ws.Range("F91").Value = "Al Gharafa"
$ws.Range("G90").Value = "AlWakrah SC"
$ws.Range("G91").Value = "Al Sadd"
$ws.Range("H90").Value = 3
$ws.Range("H91").Value = 2
$ws.Range("I90").Value = 0
$ws.Range("I91").Value = 2
$ws.Range("J90").Value = "H"
$ws.Range("J91").Value = "D"
$ws.Range("K90").Value = 2
$ws.Range("K91").Value = 5
$ws.Range("L90").Value = 3.6
$ws.Range("L91").Value = 4.75
$ws.Range("M90").Value = 3.1
$ws.Range("M91").Value = 1.45
$ws.Range("N90").Value = 2.15
$ws.Range("N91").Value = 5.25
$ws.Range("O90").Value = 3.4
$ws.Range("O91").Value = 5
$ws.Range("P90").Value = 2.9
$ws.Range("P91").Value = 1.4
$ws.Range("Q90").Value = -0.25
$ws.Range("Q91").Value = 1.25
$ws.Range("R90").Value = 1.975
$ws.Range("R91").Value = 2
$ws.Range("S90").Value = 1.825
$ws.Range("S91").Value = 1.8
$ws.Range("T90").Value = 3
$ws.Range("T91").Value = 3.75
$ws.Range("U90").Value = 1.925
$ws.Range("U91").Value = 1.875
$ws.Range("V90").Value = 1.875
$ws.Range("V91").Value = 1.925
$ws.Range("W90").Value = 1.15
$ws.Range("W91").Value = -1
$ws.Range("X90").Value = -1
$ws.Range("X91").Value = 4
$ws.Range("Y90").Value = -1
$ws.Range("Y91").Value = -1
$ws.Range("Z90").Value = 0.9750000000000001
$ws.Range("Z91").Value = 1
$ws.Range("AA90").Value = -1
$ws.Range("AA91").Value = -1
$ws.Range("AB90").Value = 0
$ws.Range("AB91").Value = 0.4375
$ws.Range("AC90").Value = -0
$ws.Range("AC91").Value = -0.5

# Swap row 92 <-> row 93
$ws.Range("B92").Value = 7840687
$ws.Range("B93").Value = 7840811
$ws.Range("F92").Value = "AlWakrah SC"
$ws.Range("F93").Value = "AlArabi Doha"
$ws.Range("G92").Value = "Al Duhail"
$ws.Range("G93").Value = "Al Sadd"
$ws.Range("H92").Value = 2
$ws.Range("H93").Value = 2
$ws.Range("I92").Value = 1
$ws.Range("I93").Value = 2
$ws.Range("J92").Value = "H"
$ws.Range("J93").Value = "D"
$ws.Range("K92").Value = 2.7
$ws.Range("K93").Value = 4.75
$ws.Range("L92").Value = 3.6
$ws.Range("L93").Value = 4.2
$ws.Range("M92").Value = 2.25
$ws.Range("M93").Value = 1.55
$ws.Range("N92").Value = 3
$ws.Range("N93").Value = 3.6
$ws.Range("O92").Value = 3.6
$ws.Range("O93").Value = 4
$ws.Range("P92").Value = 2.05
$ws.Range("P93").Value = 1.75
$ws.Range("Q92").Value = 0.25
$ws.Range("Q93").Value = 0.75
$ws.Range("R92").Value = 1.95
$ws.Range("R93").Value = 1.85
$ws.Range("S92").Value = 1.75
$ws.Range("S93").Value = 1.95
$ws.Range("T92").Value = 3.25
$ws.Range("T93").Value = 3.5
$ws.Range("U92").Value = 1.95
$ws.Range("U93").Value = 1.975
$ws.Range("V92").Value = 1.75
$ws.Range("V93").Value = 1.825
$ws.Range("W92").Value = 2
$ws.Range("W93").Value = -1
$ws.Range("X92").Value = -1
$ws.Range("X93").Value = 3
$ws.Range("Y92").Value = -1
$ws.Range("Y93").Value = -1
$ws.Range("Z92").Value = 0.95
$ws.Range("Z93").Value = 0.8500000000000001
$ws.Range("AA92").Value = -1
$ws.Range("AA93").Value = -1
$ws.Range("AB92").Value = -0.5
$ws.Range("AB93").Value = 0.9750000000000001
$ws.Range("AC92").Value = 0.375
$ws.Range("AC93").Value = -1

# Swap row 98 <-> row 99
$ws.Range("B98").Value = 7840816
$ws.Range("B99").Value = 7840688
$ws.Range("F98").Value = "Al Sadd"
$ws.Range("F99").Value = "Al Duhail"
$ws.Range("G98").Value = "AlMuaidar"
$ws.Range("G99").Value = "AlShamal SC"
$ws.Range("H98").Value = 4
$ws.Range("H99").Value = 3
$ws.Range("I98").Value = 2
$ws.Range("I99").Value = 1
$ws.Range("J98").Value = "H"
$ws.Range("J99").Value = "H"
$ws.Range("K98").Value = 1.166
$ws.Range("K99").Value = 1.444
$ws.Range("L98").Value = 7
$ws.Range("L99").Value = 4.75
$ws.Range("M98").Value = 11
$ws.Range("M99").Value = 5.5
$ws.Range("N98").Value = 1.285
$ws.Range("N99").Value = 1.615
$ws.Range("O98").Value = 5.5
$ws.Range("O99").Value = 4.2
$ws.Range("P98").Value = 8
$ws.Range("P99").Value = 4.333
$ws.Range("Q98").Value = -1.75
$ws.Range("Q99").Value = -0.75
$ws.Range("R98").Value = 1.975
$ws.Range("R99").Value = 1.8
$ws.Range("S98").Value = 1.825
$ws.Range("S99").Value = 2
$ws.Range("T98").Value = 3.75
$ws.Range("T99").Value = 3.25
$ws.Range("U98").Value = 1.925
$ws.Range("U99").Value = 1.925
$ws.Range("V98").Value = 1.775
$ws.Range("V99").Value = 1.875
$ws.Range("W98").Value = 0.2849999999999999
$ws.Range("W99").Value = 0.615
$ws.Range("X98").Value = -1
$ws.Range("X99").Value = -1
$ws.Range("Y98").Value = -1
$ws.Range("Y99").Value = -1
$ws.Range("Z98").Value = 0.4875
$ws.Range("Z99").Value = 0.8
$ws.Range("AA98").Value = -0.5
$ws.Range("AA99").Value = -1
$ws.Range("AB98").Value = 0.925
$ws.Range("AB99").Value = 0.925
$ws.Range("AC98").Value = -1
$ws.Range("AC99").Value = -1

# Swap row 100 <-> row 101
$ws.Range("B100").Value = 7003492
$ws.Range("B101").Value = 7004650
$ws.Range("F100").Value = "AlAhli Doha"
$ws.Range("F101").Value = "Umm Salal"
$ws.Range("G100").Value = "AlArabi Doha"
$ws.Range("G101").Value = "AlWakrah SC"
$ws.Range("H100").Value = 1
$ws.Range("H101").Value = 0
$ws.Range("I100").Value = 1
$ws.Range("I101").Value = 2
$ws.Range("J100").Value = "D"
$ws.Range("J101").Value = "A"
$ws.Range("K100").Value = 4.5
$ws.Range("K101").Value = 3.8
$ws.Range("L100").Value = 4.333
$ws.Range("L101").Value = 3.75
$ws.Range("M100").Value = 1.55
$ws.Range("M101").Value = 1.8
$ws.Range("N100").Value = 4
$ws.Range("N101").Value = 4.75
$ws.Range("O100").Value = 4
$ws.Range("O101").Value = 4
$ws.Range("P100").Value = 1.666
$ws.Range("P101").Value = 1.571
$ws.Range("Q100").Value = 0.75
$ws.Range("Q101").Value = 1
$ws.Range("R100").Value = 1.95
$ws.Range("R101").Value = 1.825
$ws.Range("S100").Value = 1.85
$ws.Range("S101").Value = 1.975
$ws.Range("T100").Value = 3.25
$ws.Range("T101").Value = 3
$ws.Range("U100").Value = 1.775
$ws.Range("U101").Value = 1.85
$ws.Range("V100").Value = 1.925
$ws.Range("V101").Value = 1.95
$ws.Range("W100").Value = -1
$ws.Range("W101").Value = -1
$ws.Range("X100").Value = 3
$ws.Range("X101").Value = -1
$ws.Range("Y100").Value = -1
$ws.Range("Y101").Value = 0.571
$ws.Range("Z100").Value = 0.95
$ws.Range("Z101").Value = -1
$ws.Range("AA100").Value = -1
$ws.Range("AA101").Value = 0.9750000000000001
$ws.Range("AB100").Value = -1
$ws.Range("AB101").Value = -1
$ws.Range("AC100").Value = 0.925
$ws.Range("AC101").Value = 0.95

# Swap row 102 <-> row 103
$ws.Range("B102").Value = 7840817
$ws.Range("B103").Value = 7840818
$ws.Range("F102").Value = "Al Markhiya"
$ws.Range("F103").Value = "Al Gharafa"
$ws.Range("G102").Value = "Qatar SC Doha"
$ws.Range("G103").Value = "AlRayyan SC"
$ws.Range("H102").Value = 1
$ws.Range("H103").Value = 3
$ws.Range("I102").Value = 2
$ws.Range("I103").Value = 0
$ws.Range("J102").Value = "A"
$ws.Range("J103").Value = "H"
$ws.Range("K102").Value = 4
$ws.Range("K103").Value = 2.5
$ws.Range("L102").Value = 3.6
$ws.Range("L103").Value = 3.4
$ws.Range("M102").Value = 1.75
$ws.Range("M103").Value = 2.5
$ws.Range("N102").Value = 3.5
$ws.Range("N103").Value = 2.3
$ws.Range("O102").Value = 3.75
$ws.Range("O103").Value = 3.5
$ws.Range("P102").Value = 1.833
$ws.Range("P103").Value = 2.7
$ws.Range("Q102").Value = 0.5
$ws.Range("Q103").Value = 0
$ws.Range("R102").Value = 1.975
$ws.Range("R103").Value = 1.75
$ws.Range("S102").Value = 1.825
$ws.Range("S103").Value = 2.05
$ws.Range("T102").Value = 3
$ws.Range("T103").Value = 3.5
$ws.Range("U102").Value = 1.8
$ws.Range("U103").Value = 1.925
$ws.Range("V102").Value = 2
$ws.Range("V103").Value = 1.875
$ws.Range("W102").Value = -1
$ws.Range("W103").Value = 1.3
$ws.Range("X102").Value = -1
$ws.Range("X103").Value = -1
$ws.Range("Y102").Value = 0.833
$ws.Range("Y103").Value = -1
$ws.Range("Z102").Value = -1
$ws.Range("Z103").Value = 0.75
$ws.Range("AA102").Value = 0.825
$ws.Range("AA103").Value = -1
$ws.Range("AB102").Value = 0
$ws.Range("AB103").Value = -1
$ws.Range("AC102").Value = -0
$ws.Range("AC103").Value = 0.875

# Swap row 106 <-> row 107
$ws.Range("B106").Value = 7004653
$ws.Range("B107").Value = 7004656
$ws.Range("F106").Value = "Qatar SC Doha"
$ws.Range("F107").Value = "AlArabi Doha"
$ws.Range("G106").Value = "AlWakrah SC"
$ws.Range("G107").Value = "Al Gharafa"
$ws.Range("H106").Value = 1
$ws.Range("H107").Value = 0
$ws.Range("I106").Value = 5
$ws.Range("I107").Value = 1
$ws.Range("J106").Value = "A"
$ws.Range("J107").Value = "A"
$ws.Range("K106").Value = 3.5
$ws.Range("K107").Value = 2.55
$ws.Range("L106").Value = 3.6
$ws.Range("L107").Value = 3.6
$ws.Range("M106").Value = 1.85
$ws.Range("M107").Value = 2.3
$ws.Range("N106").Value = 3.4
$ws.Range("N107").Value = 2.4
$ws.Range("O106").Value = 3.6
$ws.Range("O107").Value = 3.6
$ws.Range("P106").Value = 1.85
$ws.Range("P107").Value = 2.375
$ws.Range("Q106").Value = 0.5
$ws.Range("Q107").Value = 0
$ws.Range("R106").Value = 1.9
$ws.Range("R107").Value = 1.95
$ws.Range("S106").Value = 1.9
$ws.Range("S107").Value = 1.85
$ws.Range("T106").Value = 3
$ws.Range("T107").Value = 3.25
$ws.Range("U106").Value = 1.925
$ws.Range("U107").Value = 1.775
$ws.Range("V106").Value = 1.875
$ws.Range("V107").Value = 1.925
$ws.Range("W106").Value = -1
$ws.Range("W107").Value = -1
$ws.Range("X106").Value = -1
$ws.Range("X107").Value = -1
$ws.Range("Y106").Value = 0.8500000000000001
$ws.Range("Y107").Value = 1.375
$ws.Range("Z106").Value = -1
$ws.Range("Z107").Value = -1
$ws.Range("AA106").Value = 0.8999999999999999
$ws.Range("AA107").Value = 0.8500000000000001
$ws.Range("AB106").Value = 0.925
$ws.Range("AB107").Value = -1
$ws.Range("AC106").Value = -1
$ws.Range("AC107").Value = 0.925

# Swap row 110 <-> row 111
$ws.Range("B110").Value = 8022181
$ws.Range("B111").Value = 7004660
$ws.Range("F110").Value = "Al Duhail"
$ws.Range("F111").Value = "AlWakrah SC"
$ws.Range("G110").Value = "Al Sadd"
$ws.Range("G111").Value = "AlArabi Doha"
$ws.Range("H110").Value = 3
$ws.Range("H111").Value = 2
$ws.Range("I110").Value = 1
$ws.Range("I111").Value = 4
$ws.Range("J110").Value = "H"
$ws.Range("J111").Value = "A"
$ws.Range("K110").Value = 4.5
$ws.Range("K111").Value = 2
$ws.Range("L110").Value = 4.333
$ws.Range("L111").Value = 3.75
$ws.Range("M110").Value = 1.6
$ws.Range("M111").Value = 3.2
$ws.Range("N110").Value = 4
$ws.Range("N111").Value = 1.909
$ws.Range("O110").Value = 4.2
$ws.Range("O111").Value = 3.8
$ws.Range("P110").Value = 1.666
$ws.Range("P111").Value = 3.4
$ws.Range("Q110").Value = 0.75
$ws.Range("Q111").Value = -0.5
$ws.Range("R110").Value = 1.95
$ws.Range("R111").Value = 1.95
$ws.Range("S110").Value = 1.85
$ws.Range("S111").Value = 1.85
$ws.Range("T110").Value = 3.5
$ws.Range("T111").Value = 3.25
$ws.Range("U110").Value = 1.85
$ws.Range("U111").Value = 2
$ws.Range("V110").Value = 1.95
$ws.Range("V111").Value = 1.8
$ws.Range("W110").Value = 3
$ws.Range("W111").Value = -1
$ws.Range("X110").Value = -1
$ws.Range("X111").Value = -1
$ws.Range("Y110").Value = -1
$ws.Range("Y111").Value = 2.4
$ws.Range("Z110").Value = 0.95
$ws.Range("Z111").Value = -1
$ws.Range("AA110").Value = -1
$ws.Range("AA111").Value = 0.8500000000000001
$ws.Range("AB110").Value = 0.8500000000000001
$ws.Range("AB111").Value = 1
$ws.Range("AC110").Value = -1
$ws.Range("AC111").Value = -1

# Swap row 114 <-> row 115
$ws.Range("B114").Value = 7004658
$ws.Range("B115").Value = 7004659
$ws.Range("F114").Value = "Al Gharafa"
$ws.Range("F115").Value = "AlRayyan SC"
$ws.Range("G114").Value = "Umm Salal"
$ws.Range("G115").Value = "AlMuaidar"
$ws.Range("H114").Value = 1
$ws.Range("H115").Value = 1
$ws.Range("I114").Value = 1
$ws.Range("I115").Value = 0
$ws.Range("J114").Value = "D"
$ws.Range("J115").Value = "H"
$ws.Range("K114").Value = 1.571
$ws.Range("K115").Value = 1.5
$ws.Range("L114").Value = 3.25
$ws.Range("L115").Value = 3.6
$ws.Range("M114").Value = 6
$ws.Range("M115").Value = 6
$ws.Range("N114").Value = 1.45
$ws.Range("N115").Value = 1.5
$ws.Range("O114").Value = 3.6
$ws.Range("O115").Value = 3.8
$ws.Range("P114").Value = 6.5
$ws.Range("P115").Value = 5.75
$ws.Range("Q114").Value = -1.25
$ws.Range("Q115").Value = -1
$ws.Range("R114").Value = 1.95
$ws.Range("R115").Value = 1.8
$ws.Range("S114").Value = 1.85
$ws.Range("S115").Value = 2
$ws.Range("T114").Value = 3.25
$ws.Range("T115").Value = 3.25
$ws.Range("U114").Value = 1.85
$ws.Range("U115").Value = 1.9
$ws.Range("V114").Value = 1.95
$ws.Range("V115").Value = 1.9
$ws.Range("W114").Value = -1
$ws.Range("W115").Value = 0.5
$ws.Range("X114").Value = 2.6
$ws.Range("X115").Value = -1
$ws.Range("Y114").Value = -1
$ws.Range("Y115").Value = -1
$ws.Range("Z114").Value = -1
$ws.Range("Z115").Value = 0
$ws.Range("AA114").Value = 0.8500000000000001
$ws.Range("AA115").Value = -0
$ws.Range("AB114").Value = -1
$ws.Range("AB115").Value = -1
$ws.Range("AC114").Value = 0.95
$ws.Range("AC115").Value = 0.8999999999999999
